$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.743.34"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").Value = "'1.851.67"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("D4").Value = "'1.036"
$ws.Range("E4").Value = "  +0.89%  "

$ws.Range("D5").Value = "'322.80"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("D6").Value = "'1.032"
$ws.Range("E6").Value = "  +0.65%  "

$ws.Range("D7").Value = "'0.4396"
$ws.Range("E7").Value = "  +0.54%  "

$ws.Range("D8").Value = "'0.3795"
$ws.Range("E8").Value = "  +1.54%  "

$ws.Range("D9").Value = "'0.07394"
$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").Value = "'0.8833"
$ws.Range("E10").Value = "  +0.96%  "

$ws.Range("D11").Value = "'21.55"
$ws.Range("E11").Value = "  +0.21%  "

$ws.Range("D12").Value = "'1.858.17"
$ws.Range("E12").Value = "  -0.04%  "

$ws.Range("D13").Value = "'5.500"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("D14").Value = "'6.704"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").Value = "'0.07173"
$ws.Range("E15").Value = "  +0.30%  "

$ws.Range("D16").Value = "'85.02"
$ws.Range("E16").Value = "  +2.83%  "

$ws.Range("D17").Value = "'1.039"
$ws.Range("E17").Value = "  +0.87%  "

$ws.Range("D18").Value = "'0.000009070"
$ws.Range("E18").Value = "  +0.77%  "

$ws.Range("D19").Value = "'1.033"
$ws.Range("E19").Value = "  +0.70%  "

$ws.Range("D20").Value = "'15.47"
$ws.Range("E20").Value = "  +0.45%  "

$ws.Range("D21").Value = "'27.745.15"
$ws.Range("E21").Value = "  +0.76%  "

$ws.Range("E22").Value = "  +0.60%  "

$ws.Range("E23").Value = "  +1.06%  "

$ws.Range("D24").Value = "'2.080.97"
$ws.Range("E24").Value = "  +1.08%  "

$ws.Range("D25").Value = "'2.058"
$ws.Range("E25").Value = "  +6.80%  "

$ws.Range("D26").Value = "'158.90"
$ws.Range("E26").Value = "  +0.91%  "

$ws.Range("D27").Value = "'18.70"
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").Value = "'1.992"
$ws.Range("E28").Value = "  +2.83%  "

$ws.Range("D29").Value = "'5.322"
$ws.Range("E29").Value = "  +1.28%  "

$ws.Range("D30").Value = "'117.81"
$ws.Range("E30").Value = "  +1.65%  "

$ws.Range("D31").Value = "'0.09090"
$ws.Range("E31").Value = "  -0.14%  "

$ws.Range("D32").Value = "'0.7723"
$ws.Range("E32").Value = "  +0.62%  "

$ws.Range("D33").Value = "'1.210"
$ws.Range("E33").Value = "  +0.30%  "

$ws.Range("D34").Value = "'3.010"
$ws.Range("E34").Value = "  +4.84%  "

$ws.Range("D35").Value = "'4.555"
$ws.Range("E35").Value = "  +1.17%  "

$ws.Range("D36").Value = "'1.034"
$ws.Range("E36").Value = "  +0.65%  "

$ws.Range("D37").Value = "'1.150"
$ws.Range("E37").Value = "  +0.74%  "

$ws.Range("D38").Value = "'0.01971"
$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("D39").Value = "'0.05261"
$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("D40").Value = "'2.850"
$ws.Range("E40").Value = "  +2.52%  "

$ws.Range("D41").Value = "'0.5172"
$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("D42").Value = "'0.1669"
$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("D43").Value = "'6.875"
$ws.Range("E43").Value = "  +3.29%  "

$ws.Range("D44").Value = "'8.711"
$ws.Range("E44").Value = "  +2.09%  "

$ws.Range("D45").Value = "'110.26"
$ws.Range("E45").Value = "  +1.11%  "

$ws.Range("D46").Value = "'10.75"
$ws.Range("E46").Value = "  +1.97%  "

$ws.Range("D47").Value = "'1.036"
$ws.Range("E47").Value = "  +0.69%  "

$ws.Range("D48").Value = "'0.06568"
$ws.Range("E48").Value = "  +3.19%  "

$ws.Range("D49").Value = "'1.702"
$ws.Range("E49").Value = "  -0.61%  "

$ws.Range("D50").Value = "'0.4696"
$ws.Range("E50").Value = "  +1.10%  "

$ws.Range("D51").Value = "'1.885"
$ws.Range("E51").Value = "  -0.58%  "

